# Submit and email excel feature is fully functional
# Update product codes on the Order sheet and add a new line item.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 keeps quantity 5, but the product code text changes.
$ws.Range("B2").Value = "DS2310BLK-LF"

# Row 3 switches to the other product code (quantity 10 unchanged).
$ws.Range("B3").Value = "DS2310WMUS-LF"

# New row 4: same product code as row 3, quantity 4.
$ws.Range("B4").Value = "DS2310WMUS-LF"
$ws.Range("C4").Value = 4
